$d = $word.ActiveDocument

# Locate the unique anchor text that spans the region we need to reformat.
# The original run's text is:
#   " is simply dummy text of the printing and typesetting industry. Lorem
#    Ipsum has been the industry's standard dummy text ever since the 1500s,
#    when an unknown printer took a galley of "
# We only need to change the character formatting of two sub-spans:
#   "dummy text" (2nd occurrence)  -> sz 30 (15pt), darkCyan highlight
#   "sin" (inside "since")          -> darkMagenta highlight
# No characters are inserted/deleted/replaced; only run-splitting + direct
# character formatting is applied, so the document's plain text is unchanged.

$anchorRange = $d.Content
$found = $anchorRange.Find.Execute("standard dummy text ever since the 1500s")
if (-not $found) {
    throw "Could not locate anchor text for formatting the Lorem Ipsum paragraph."
}

$anchorStart = $anchorRange.Start

# Offsets of the two sub-spans to format, relative to $anchorStart
# (computed against the literal anchor string above).
$dummyTextStart = $anchorStart + 9
$dummyTextEnd   = $dummyTextStart + 10   # "dummy text" = 10 chars

$sinStart = $anchorStart + 25
$sinEnd   = $sinStart + 3                 # "sin" = 3 chars

# --- Format "dummy text" (2nd occurrence): 15pt, dark cyan highlight ---
$dummyTextRange = $d.Range($dummyTextStart, $dummyTextEnd)
if ($dummyTextRange.Text -ne "dummy text") {
    throw "Unexpected text at dummy-text range: [$($dummyTextRange.Text)]"
}
$dummyTextRange.Font.Size = 15
$dummyTextRange.Font.HighlightColorIndex = "darkCyan"

# --- Format "sin" (inside "since"): dark magenta highlight ---
$sinRange = $d.Range($sinStart, $sinEnd)
if ($sinRange.Text -ne "sin") {
    throw "Unexpected text at sin range: [$($sinRange.Text)]"
}
$sinRange.Font.HighlightColorIndex = "darkMagenta"
